$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1): new category values, drop columns G:L ---
$ws.Cells.Item(1,1).Value = "p"
$ws.Cells.Item(1,2).Value = "b"
$ws.Cells.Item(1,3).Value = "Y"
$ws.Cells.Item(1,4).Value = "f"
$ws.Cells.Item(1,5).Value = "f"
$ws.Cells.Item(1,6).Value = " d"
$ws.Cells.Item(1,13).Value = "s"
$ws.Cells.Item(1,14).Value = "s"
$ws.Cells.Item(1,15).Value = "w"
$ws.Cells.Item(1,16).Value = "w"
$ws.Cells.Item(1,17).Value = "p"
$ws.Cells.Item(1,18).Value = "w"
$ws.Cells.Item(1,19).Value = "o"
$ws.Cells.Item(1,20).Value = "p"
$ws.Cells.Item(1,21).Value = "k"
$ws.Cells.Item(1,22).Value = "s"
$ws.Cells.Item(1,23).Value = "u"
$ws.Range("G1:L1").ClearContents()

# --- Append new rows 17-32 (row17 = new header values repeated; rows 18-32 = old rows 2-16 data) ---
# Row 17
$ws.Cells.Item(17,1).Value = "p"
$ws.Cells.Item(17,2).Value = "b"
$ws.Cells.Item(17,3).Value = "Y"
$ws.Cells.Item(17,4).Value = "f"
$ws.Cells.Item(17,5).Value = "f"
$ws.Cells.Item(17,6).Value = " d"
$ws.Cells.Item(17,13).Value = "s"
$ws.Cells.Item(17,14).Value = "s"
$ws.Cells.Item(17,15).Value = "w"
$ws.Cells.Item(17,16).Value = "w"
$ws.Cells.Item(17,17).Value = "p"
$ws.Cells.Item(17,18).Value = "w"
$ws.Cells.Item(17,19).Value = "o"
$ws.Cells.Item(17,20).Value = "p"
$ws.Cells.Item(17,21).Value = "k"
$ws.Cells.Item(17,22).Value = "s"
$ws.Cells.Item(17,23).Value = "u"

# Row 18
$ws.Cells.Item(18,1).Value = "e"
$ws.Cells.Item(18,2).Value = "x"
$ws.Cells.Item(18,3).Value = "s"
$ws.Cells.Item(18,4).Value = "y"
$ws.Cells.Item(18,5).Value = "t"
$ws.Cells.Item(18,6).Value = "a"
$ws.Cells.Item(18,7).Value = "f"
$ws.Cells.Item(18,8).Value = "c"
$ws.Cells.Item(18,9).Value = "b"
$ws.Cells.Item(18,10).Value = "k"
$ws.Cells.Item(18,11).Value = "e"
$ws.Cells.Item(18,12).Value = "c"
$ws.Cells.Item(18,13).Value = "s"
$ws.Cells.Item(18,14).Value = "s"
$ws.Cells.Item(18,15).Value = "w"
$ws.Cells.Item(18,16).Value = "w"
$ws.Cells.Item(18,17).Value = "p"
$ws.Cells.Item(18,18).Value = "w"
$ws.Cells.Item(18,19).Value = "o"
$ws.Cells.Item(18,20).Value = "p"
$ws.Cells.Item(18,21).Value = "n"
$ws.Cells.Item(18,22).Value = "n"
$ws.Cells.Item(18,23).Value = "g"

# Row 19
$ws.Cells.Item(19,1).Value = "e"
$ws.Cells.Item(19,2).Value = "b"
$ws.Cells.Item(19,3).Value = "s"
$ws.Cells.Item(19,4).Value = "w"
$ws.Cells.Item(19,5).Value = "t"
$ws.Cells.Item(19,6).Value = "l"
$ws.Cells.Item(19,7).Value = "f"
$ws.Cells.Item(19,8).Value = "c"
$ws.Cells.Item(19,9).Value = "b"
$ws.Cells.Item(19,10).Value = "n"
$ws.Cells.Item(19,11).Value = "e"
$ws.Cells.Item(19,12).Value = "c"
$ws.Cells.Item(19,13).Value = "s"
$ws.Cells.Item(19,14).Value = "s"
$ws.Cells.Item(19,15).Value = "w"
$ws.Cells.Item(19,16).Value = "w"
$ws.Cells.Item(19,17).Value = "p"
$ws.Cells.Item(19,18).Value = "w"
$ws.Cells.Item(19,19).Value = "o"
$ws.Cells.Item(19,20).Value = "p"
$ws.Cells.Item(19,21).Value = "n"
$ws.Cells.Item(19,22).Value = "n"
$ws.Cells.Item(19,23).Value = "m"

# Row 20
$ws.Cells.Item(20,1).Value = "p"
$ws.Cells.Item(20,2).Value = "x"
$ws.Cells.Item(20,3).Value = "y"
$ws.Cells.Item(20,4).Value = "w"
$ws.Cells.Item(20,5).Value = "t"
$ws.Cells.Item(20,6).Value = "p"
$ws.Cells.Item(20,7).Value = "f"
$ws.Cells.Item(20,8).Value = "c"
$ws.Cells.Item(20,9).Value = "n"
$ws.Cells.Item(20,10).Value = "n"
$ws.Cells.Item(20,11).Value = "e"
$ws.Cells.Item(20,12).Value = "e"
$ws.Cells.Item(20,13).Value = "s"
$ws.Cells.Item(20,14).Value = "s"
$ws.Cells.Item(20,15).Value = "w"
$ws.Cells.Item(20,16).Value = "w"
$ws.Cells.Item(20,17).Value = "p"
$ws.Cells.Item(20,18).Value = "w"
$ws.Cells.Item(20,19).Value = "o"
$ws.Cells.Item(20,20).Value = "p"
$ws.Cells.Item(20,21).Value = "k"
$ws.Cells.Item(20,22).Value = "s"
$ws.Cells.Item(20,23).Value = "u"

# Row 21
$ws.Cells.Item(21,1).Value = "e"
$ws.Cells.Item(21,2).Value = "x"
$ws.Cells.Item(21,3).Value = "s"
$ws.Cells.Item(21,4).Value = "g"
$ws.Cells.Item(21,5).Value = "f"
$ws.Cells.Item(21,6).Value = "n"
$ws.Cells.Item(21,7).Value = "f"
$ws.Cells.Item(21,8).Value = "w"
$ws.Cells.Item(21,9).Value = "b"
$ws.Cells.Item(21,10).Value = "k"
$ws.Cells.Item(21,11).Value = "t"
$ws.Cells.Item(21,12).Value = "e"
$ws.Cells.Item(21,13).Value = "s"
$ws.Cells.Item(21,14).Value = "s"
$ws.Cells.Item(21,15).Value = "w"
$ws.Cells.Item(21,16).Value = "w"
$ws.Cells.Item(21,17).Value = "p"
$ws.Cells.Item(21,18).Value = "w"
$ws.Cells.Item(21,19).Value = "o"
$ws.Cells.Item(21,20).Value = "e"
$ws.Cells.Item(21,21).Value = "n"
$ws.Cells.Item(21,22).Value = "a"
$ws.Cells.Item(21,23).Value = "g"

# Row 22
$ws.Cells.Item(22,1).Value = "e"
$ws.Cells.Item(22,2).Value = "x"
$ws.Cells.Item(22,3).Value = "y"
$ws.Cells.Item(22,4).Value = "y"
$ws.Cells.Item(22,5).Value = "t"
$ws.Cells.Item(22,6).Value = "a"
$ws.Cells.Item(22,7).Value = "f"
$ws.Cells.Item(22,8).Value = "c"
$ws.Cells.Item(22,9).Value = "b"
$ws.Cells.Item(22,10).Value = "n"
$ws.Cells.Item(22,11).Value = "e"
$ws.Cells.Item(22,12).Value = "c"
$ws.Cells.Item(22,13).Value = "s"
$ws.Cells.Item(22,14).Value = "s"
$ws.Cells.Item(22,15).Value = "w"
$ws.Cells.Item(22,16).Value = "w"
$ws.Cells.Item(22,17).Value = "p"
$ws.Cells.Item(22,18).Value = "w"
$ws.Cells.Item(22,19).Value = "o"
$ws.Cells.Item(22,20).Value = "p"
$ws.Cells.Item(22,21).Value = "k"
$ws.Cells.Item(22,22).Value = "n"
$ws.Cells.Item(22,23).Value = "g"

# Row 23
$ws.Cells.Item(23,1).Value = "e"
$ws.Cells.Item(23,2).Value = "b"
$ws.Cells.Item(23,3).Value = "s"
$ws.Cells.Item(23,4).Value = "w"
$ws.Cells.Item(23,5).Value = "t"
$ws.Cells.Item(23,6).Value = "a"
$ws.Cells.Item(23,7).Value = "f"
$ws.Cells.Item(23,8).Value = "c"
$ws.Cells.Item(23,9).Value = "b"
$ws.Cells.Item(23,10).Value = "g"
$ws.Cells.Item(23,11).Value = "e"
$ws.Cells.Item(23,12).Value = "c"
$ws.Cells.Item(23,13).Value = "s"
$ws.Cells.Item(23,14).Value = "s"
$ws.Cells.Item(23,15).Value = "w"
$ws.Cells.Item(23,16).Value = "w"
$ws.Cells.Item(23,17).Value = "p"
$ws.Cells.Item(23,18).Value = "w"
$ws.Cells.Item(23,19).Value = "o"
$ws.Cells.Item(23,20).Value = "p"
$ws.Cells.Item(23,21).Value = "k"
$ws.Cells.Item(23,22).Value = "n"
$ws.Cells.Item(23,23).Value = "m"

# Row 24
$ws.Cells.Item(24,1).Value = "e"
$ws.Cells.Item(24,2).Value = "b"
$ws.Cells.Item(24,3).Value = "y"
$ws.Cells.Item(24,4).Value = "w"
$ws.Cells.Item(24,5).Value = "t"
$ws.Cells.Item(24,6).Value = "l"
$ws.Cells.Item(24,7).Value = "f"
$ws.Cells.Item(24,8).Value = "c"
$ws.Cells.Item(24,9).Value = "b"
$ws.Cells.Item(24,10).Value = "n"
$ws.Cells.Item(24,11).Value = "e"
$ws.Cells.Item(24,12).Value = "c"
$ws.Cells.Item(24,13).Value = "s"
$ws.Cells.Item(24,14).Value = "s"
$ws.Cells.Item(24,15).Value = "w"
$ws.Cells.Item(24,16).Value = "w"
$ws.Cells.Item(24,17).Value = "p"
$ws.Cells.Item(24,18).Value = "w"
$ws.Cells.Item(24,19).Value = "o"
$ws.Cells.Item(24,20).Value = "p"
$ws.Cells.Item(24,21).Value = "n"
$ws.Cells.Item(24,22).Value = "s"
$ws.Cells.Item(24,23).Value = "m"

# Row 25
$ws.Cells.Item(25,1).Value = "p"
$ws.Cells.Item(25,2).Value = "x"
$ws.Cells.Item(25,3).Value = "y"
$ws.Cells.Item(25,4).Value = "w"
$ws.Cells.Item(25,5).Value = "t"
$ws.Cells.Item(25,6).Value = "p"
$ws.Cells.Item(25,7).Value = "f"
$ws.Cells.Item(25,8).Value = "c"
$ws.Cells.Item(25,9).Value = "n"
$ws.Cells.Item(25,10).Value = "p"
$ws.Cells.Item(25,11).Value = "e"
$ws.Cells.Item(25,12).Value = "e"
$ws.Cells.Item(25,13).Value = "s"
$ws.Cells.Item(25,14).Value = "s"
$ws.Cells.Item(25,15).Value = "w"
$ws.Cells.Item(25,16).Value = "w"
$ws.Cells.Item(25,17).Value = "p"
$ws.Cells.Item(25,18).Value = "w"
$ws.Cells.Item(25,19).Value = "o"
$ws.Cells.Item(25,20).Value = "p"
$ws.Cells.Item(25,21).Value = "k"
$ws.Cells.Item(25,22).Value = "v"
$ws.Cells.Item(25,23).Value = "g"

# Row 26
$ws.Cells.Item(26,1).Value = "e"
$ws.Cells.Item(26,2).Value = "b"
$ws.Cells.Item(26,3).Value = "s"
$ws.Cells.Item(26,4).Value = "y"
$ws.Cells.Item(26,5).Value = "t"
$ws.Cells.Item(26,6).Value = "a"
$ws.Cells.Item(26,7).Value = "f"
$ws.Cells.Item(26,8).Value = "c"
$ws.Cells.Item(26,9).Value = "b"
$ws.Cells.Item(26,10).Value = "g"
$ws.Cells.Item(26,11).Value = "e"
$ws.Cells.Item(26,12).Value = "c"
$ws.Cells.Item(26,13).Value = "s"
$ws.Cells.Item(26,14).Value = "s"
$ws.Cells.Item(26,15).Value = "w"
$ws.Cells.Item(26,16).Value = "w"
$ws.Cells.Item(26,17).Value = "p"
$ws.Cells.Item(26,18).Value = "w"
$ws.Cells.Item(26,19).Value = "o"
$ws.Cells.Item(26,20).Value = "p"
$ws.Cells.Item(26,21).Value = "k"
$ws.Cells.Item(26,22).Value = "s"
$ws.Cells.Item(26,23).Value = "m"

# Row 27
$ws.Cells.Item(27,1).Value = "e"
$ws.Cells.Item(27,2).Value = "x"
$ws.Cells.Item(27,3).Value = "y"
$ws.Cells.Item(27,4).Value = "y"
$ws.Cells.Item(27,5).Value = "t"
$ws.Cells.Item(27,6).Value = "l"
$ws.Cells.Item(27,7).Value = "f"
$ws.Cells.Item(27,8).Value = "c"
$ws.Cells.Item(27,9).Value = "b"
$ws.Cells.Item(27,10).Value = "g"
$ws.Cells.Item(27,11).Value = "e"
$ws.Cells.Item(27,12).Value = "c"
$ws.Cells.Item(27,13).Value = "s"
$ws.Cells.Item(27,14).Value = "s"
$ws.Cells.Item(27,15).Value = "w"
$ws.Cells.Item(27,16).Value = "w"
$ws.Cells.Item(27,17).Value = "p"
$ws.Cells.Item(27,18).Value = "w"
$ws.Cells.Item(27,19).Value = "o"
$ws.Cells.Item(27,20).Value = "p"
$ws.Cells.Item(27,21).Value = "n"
$ws.Cells.Item(27,22).Value = "n"
$ws.Cells.Item(27,23).Value = "g"

# Row 28
$ws.Cells.Item(28,1).Value = "e"
$ws.Cells.Item(28,2).Value = "x"
$ws.Cells.Item(28,3).Value = "y"
$ws.Cells.Item(28,4).Value = "y"
$ws.Cells.Item(28,5).Value = "t"
$ws.Cells.Item(28,6).Value = "a"
$ws.Cells.Item(28,7).Value = "f"
$ws.Cells.Item(28,8).Value = "c"
$ws.Cells.Item(28,9).Value = "b"
$ws.Cells.Item(28,10).Value = "n"
$ws.Cells.Item(28,11).Value = "e"
$ws.Cells.Item(28,12).Value = "c"
$ws.Cells.Item(28,13).Value = "s"
$ws.Cells.Item(28,14).Value = "s"
$ws.Cells.Item(28,15).Value = "w"
$ws.Cells.Item(28,16).Value = "w"
$ws.Cells.Item(28,17).Value = "p"
$ws.Cells.Item(28,18).Value = "w"
$ws.Cells.Item(28,19).Value = "o"
$ws.Cells.Item(28,20).Value = "p"
$ws.Cells.Item(28,21).Value = "k"
$ws.Cells.Item(28,22).Value = "s"
$ws.Cells.Item(28,23).Value = "m"

# Row 29
$ws.Cells.Item(29,1).Value = "e"
$ws.Cells.Item(29,2).Value = "b"
$ws.Cells.Item(29,3).Value = "s"
$ws.Cells.Item(29,4).Value = "y"
$ws.Cells.Item(29,5).Value = "t"
$ws.Cells.Item(29,6).Value = "a"
$ws.Cells.Item(29,7).Value = "f"
$ws.Cells.Item(29,8).Value = "c"
$ws.Cells.Item(29,9).Value = "b"
$ws.Cells.Item(29,10).Value = "w"
$ws.Cells.Item(29,11).Value = "e"
$ws.Cells.Item(29,12).Value = "c"
$ws.Cells.Item(29,13).Value = "s"
$ws.Cells.Item(29,14).Value = "s"
$ws.Cells.Item(29,15).Value = "w"
$ws.Cells.Item(29,16).Value = "w"
$ws.Cells.Item(29,17).Value = "p"
$ws.Cells.Item(29,18).Value = "w"
$ws.Cells.Item(29,19).Value = "o"
$ws.Cells.Item(29,20).Value = "p"
$ws.Cells.Item(29,21).Value = "n"
$ws.Cells.Item(29,22).Value = "s"
$ws.Cells.Item(29,23).Value = "g"

# Row 30
$ws.Cells.Item(30,1).Value = "p"
$ws.Cells.Item(30,2).Value = "x"
$ws.Cells.Item(30,3).Value = "y"
$ws.Cells.Item(30,4).Value = "w"
$ws.Cells.Item(30,5).Value = "t"
$ws.Cells.Item(30,6).Value = "p"
$ws.Cells.Item(30,7).Value = "f"
$ws.Cells.Item(30,8).Value = "c"
$ws.Cells.Item(30,9).Value = "n"
$ws.Cells.Item(30,10).Value = "k"
$ws.Cells.Item(30,11).Value = "e"
$ws.Cells.Item(30,12).Value = "e"
$ws.Cells.Item(30,13).Value = "s"
$ws.Cells.Item(30,14).Value = "s"
$ws.Cells.Item(30,15).Value = "w"
$ws.Cells.Item(30,16).Value = "w"
$ws.Cells.Item(30,17).Value = "p"
$ws.Cells.Item(30,18).Value = "w"
$ws.Cells.Item(30,19).Value = "o"
$ws.Cells.Item(30,20).Value = "p"
$ws.Cells.Item(30,21).Value = "n"
$ws.Cells.Item(30,22).Value = "v"
$ws.Cells.Item(30,23).Value = "u"

# Row 31
$ws.Cells.Item(31,1).Value = "e"
$ws.Cells.Item(31,2).Value = "x"
$ws.Cells.Item(31,3).Value = "f"
$ws.Cells.Item(31,4).Value = "n"
$ws.Cells.Item(31,5).Value = "f"
$ws.Cells.Item(31,6).Value = "n"
$ws.Cells.Item(31,7).Value = "f"
$ws.Cells.Item(31,8).Value = "w"
$ws.Cells.Item(31,9).Value = "b"
$ws.Cells.Item(31,10).Value = "n"
$ws.Cells.Item(31,11).Value = "t"
$ws.Cells.Item(31,12).Value = "e"
$ws.Cells.Item(31,13).Value = "s"
$ws.Cells.Item(31,14).Value = "f"
$ws.Cells.Item(31,15).Value = "w"
$ws.Cells.Item(31,16).Value = "w"
$ws.Cells.Item(31,17).Value = "p"
$ws.Cells.Item(31,18).Value = "w"
$ws.Cells.Item(31,19).Value = "o"
$ws.Cells.Item(31,20).Value = "e"
$ws.Cells.Item(31,21).Value = "k"
$ws.Cells.Item(31,22).Value = "a"
$ws.Cells.Item(31,23).Value = "g"

# Row 32
$ws.Cells.Item(32,1).Value = "e"
$ws.Cells.Item(32,2).Value = "s"
$ws.Cells.Item(32,3).Value = "f"
$ws.Cells.Item(32,4).Value = "g"
$ws.Cells.Item(32,5).Value = "f"
$ws.Cells.Item(32,6).Value = "n"
$ws.Cells.Item(32,7).Value = "f"
$ws.Cells.Item(32,8).Value = "c"
$ws.Cells.Item(32,9).Value = "n"
$ws.Cells.Item(32,10).Value = "k"
$ws.Cells.Item(32,11).Value = "e"
$ws.Cells.Item(32,12).Value = "e"
$ws.Cells.Item(32,13).Value = "s"
$ws.Cells.Item(32,14).Value = "s"
$ws.Cells.Item(32,15).Value = "w"
$ws.Cells.Item(32,16).Value = "w"
$ws.Cells.Item(32,17).Value = "p"
$ws.Cells.Item(32,18).Value = "w"
$ws.Cells.Item(32,19).Value = "o"
$ws.Cells.Item(32,20).Value = "p"
$ws.Cells.Item(32,21).Value = "n"
$ws.Cells.Item(32,22).Value = "y"
$ws.Cells.Item(32,23).Value = "u"

# --- Update view/selection to match the edited workbook ---
$ws.Range("A17:W32").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
